# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.532.74"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.377.71"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'406.74"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'134.47"
$ws.Range("E6").Value = "  +7.57%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("E10").Value = "  -6.82%  "
$ws.Range("D11").Value = "'42.62"
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "3.896.49"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "'19.72"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "3.390.07"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "61.457.75"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "'11.01"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("E20").Value = "  -4.56%  "
$ws.Range("E21").Value = "  -4.58%  "
$ws.Range("D22").Value = "'85.15"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").Value = "'314.90"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'12.84"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "'4.78"
$ws.Range("E26").Value = "  +11.48%  "
$ws.Range("D27").Value = "'8.39"
$ws.Range("E27").Value = "  +6.03%  "
$ws.Range("D28").Value = "'29.59"
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "'2.60"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'11.36"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D35").Value = "'40.70"
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").Value = "'0.0483"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'51.89"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'138.91"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "'0.297"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("D46").Value = "'16.74"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'21.28"
$ws.Range("E48").Value = "  -5.50%  "
$ws.Range("D49").Value = "2.123.62"
$ws.Range("E49").Value = "  -3.97%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").Value = "'1.92"
$ws.Range("E51").Value = "  +1.32%  "
